$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last staff row (row 7, the "Test" / test@gmail.com entry) is cleared out
# as part of fixing/validating the Staff form - the row becomes blank again
# (ready for new data entry) but the hyperlink-styled C7 cell keeps its
# formatting.
$ws.Range("A7:G7").ClearContents()

# Remove the now-stale mailto hyperlink that pointed at the removed test row.
$ws.Hyperlinks.Delete()

# Update the saved cursor/selection position.
$ws.Range("C12").Select() | Out-Null
